# Update replication table (table5_panel2): refresh bootstrapped standard
# errors for theta_se / lambda_se rows and the multiple-imputation total_dof.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# theta_se row (row 4)
$ws.Range("B4").Value = "(0.48)"
$ws.Range("C4").Value = "(0.08)"

# lambda_se row (row 6)
$ws.Range("B6").Value = "(0.76)"
$ws.Range("C6").Value = "(0.04)"

# total_dof row (row 7), multiple_imputation column
$ws.Range("C7").Value = 7310
